$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsCodebook = $wb.Worksheets.Item("Codebook")

# --- Data sheet edits ---
# A4 changes from number 60 to text "sixty"
$wsData.Range("A4").Value = "sixty"

# B11 gets a new value 7000
$wsData.Range("B11").Value = 7000

# --- Codebook sheet edits ---
# Add two new rows describing additional variables (filled column by column,
# top-to-bottom within each column, to mirror the original authoring order)
$wsCodebook.Range("A5").Value = "Shoe Size"
$wsCodebook.Range("A6").Value = "Hair color"

$wsCodebook.Range("B5").Value = "Shoe Size in MEN"
$wsCodebook.Range("B6").Value = "Color of hair as of 2024/01"

$wsCodebook.Range("C5").Value = "0-24"
$wsCodebook.Range("C6").Value = "Any Color"

# --- Selections / active sheet ---
$wsData.Range("D5").Select()
$wsCodebook.Select()
$wsCodebook.Range("C10").Select()
